$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sample data row (row 2) carried a hyperlink (and its "Hyperlink" cell
# style) on E2; drop those before the row that hosts them goes away.
$ws.Range("E2").Hyperlinks.Delete()
$wb.Styles.Item("Hyperlink").Delete()

# Remove the sample data row entirely - this sheet should only carry headers.
$ws.Rows.Item(2).Delete()

# Insert a new leading column for the UserId field.
$ws.Columns.Item(1).Insert()

# Populate the header row: UserId, FirstName, LastName, UserName, Password,
# Email, UserRole, CreatedDate.
$ws.Range("A1").Value = "UserId"
$ws.Range("H1").Value = "CreatedDate"

$ws.Range("K6").Select()
